$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = "245.82"
    3 = "21.98"
    4 = "5.379"
    5 = "0.05867"
    6 = "3.389"
    7 = "6.371"
    8 = "0.8152"
    9 = "0.9977"
    10 = "0.1418"
    11 = "0.03907"
    12 = "0.07425"
    13 = "0.03041"
    14 = "4.156"
    15 = "0.09399"
    16 = "0.001600"
    17 = "0.04828"
    18 = "0.0005883"
    19 = "0.005925"
    20 = "0.004095"
    21 = "0.0009911"
    23 = "3.715"
    24 = "2.229"
    26 = "0.1294"
    27 = "0.0002490"
    40 = "0.03873"
    41 = "0.006441"
    42 = "0.1075"
    43 = "0.002597"
    44 = "0.006683"
    45 = "0.00005615"
    47 = "0.6492"
    48 = "0.1424"
    49 = "0.00002097"
}

foreach ($row in $updates.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$row]
    $cell.Style = "Normal"
}
